{"js": "// This script replaces the block of paragraphs that runs from\n// \"Pegar a descri\u00e7\u00e3o com o professor;\" through the paragraph ending in\n// \"...(Pegar quais s\u00e3o os Amp-OPs).\" with the revised content:\n//   - the old \"TODO\" paragraph is emptied out (keeping only the\n//     `_GoBack` bookmark) and restyled to match the big centered title\n//   - \"Problem\u00e1tica:\" / the problem-statement paragraph / \"Projeto:\" are\n//     each shifted up by one paragraph (taking over the style+position of\n//     the paragraph that used to precede them)\n//   - the closing \"Projeto\" paragraph is rewritten: \" testar\" is inserted\n//     after \"capaz de\", the stray mid-word bookmark around\n//     \"fu|ncionalidade\" is removed, and the closing sentence is reworded\n//     from \"Este ser\u00e1 feito ... (Pegar quais s\u00e3o os Amp-OPs).\" to\n//     \"Ele ser\u00e1 feito ... \u2013 campus Florian\u00f3polis: LM324, TL082.\"\n//\n// We locate the affected paragraphs by searching for distinctive text\n// anchors (robust against absolute paragraph-index drift) and then\n// replace the whole span in one shot via insertOoxml so the resulting\n// run/paragraph structure matches the target exactly.\n\nconst startResults = context.document.body.search(\"Pegar a descri\u00e7\u00e3o com o professor\", { matchCase: true });\nstartResults.load(\"items\");\nconst endResults = context.document.body.search(\"Pegar quais s\u00e3o os\", { matchCase: true });\nendResults.load(\"items\");\nawait context.sync();\n\nif (startResults.items.length === 0 || endResults.items.length === 0) {\n  throw new Error(\"Could not locate the paragraphs to replace.\");\n}\n\nconst startParagraph = startResults.items[0].paragraphs.getFirst();\nconst endParagraph = endResults.items[0].paragraphs.getLast();\n\nconst targetRange = startParagraph.getRange(\"Start\").expandTo(endParagraph.getRange(\"End\"));\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:jc w:val=\"center\"/><w:rPr><w:b/><w:sz w:val=\"36\"/><w:szCs w:val=\"36\"/></w:rPr></w:pPr><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p><w:p><w:pPr><w:jc w:val=\"both\"/><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>Problem\u00e1tica:</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val=\"both\"/><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>N\u00e3o \u00e9 incomum a ocorr\u00eancia de erros em aulas de laborat\u00f3rio</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> que envolvem amplificadores operacionais</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\">. Em tais situa\u00e7\u00f5es, o aluno/professor \u00e9 muitas vezes incapaz de percebe se o problema \u00e9 do </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>Amp</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>-OP utilizado</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\">, que acaba por estragar com facilidade durante as aulas pr\u00e1ticas, </w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\">ou da montagem do circuito feito pelo aluno. </w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>Na tentativa de compreender onde est\u00e1 o erro, gasta-se muito do tempo da aula.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val=\"both\"/><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>Projeto:</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val=\"both\"/><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>Este projeto visa desenvolver e montar um hardware capaz de</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> testar</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\">a </w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>funcionalidade</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> de um amplificador operacional, para a facilita\u00e7\u00e3o de aulas pr\u00e1ticas</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\">. </w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>Ele s</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>er\u00e1 feito para funcionar com os modelos utilizados no IFSC</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> \u2013 campus Florian\u00f3polis: </w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>LM324, TL082.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>`;\n\ntargetRange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# This script replaces the block of paragraphs that runs from\n# \"Pegar a descri\u00e7\u00e3o com o professor;\" through the paragraph ending in\n# \"...(Pegar quais s\u00e3o os Amp-OPs).\" with the revised content:\n#   - the old \"TODO\" paragraph is emptied out (keeping only the\n#     `_GoBack` bookmark) and restyled to match the big centered title\n#   - \"Problem\u00e1tica:\" / the problem-statement paragraph / \"Projeto:\" are\n#     each shifted up by one paragraph (taking over the style+position of\n#     the paragraph that used to precede them)\n#   - the closing \"Projeto\" paragraph is rewritten: \" testar\" is inserted\n#     after \"capaz de\", the stray mid-word bookmark around\n#     \"fu|ncionalidade\" is removed, and the closing sentence is reworded\n#     from \"Este ser\u00e1 feito ... (Pegar quais s\u00e3o os Amp-OPs).\" to\n#     \"Ele ser\u00e1 feito ... \u2013 campus Florian\u00f3polis: LM324, TL082.\"\n#\n# The affected paragraphs are located via Find (robust against absolute\n# paragraph-index drift) and then the whole span is replaced in one shot\n# via Range.InsertXML so the resulting run/paragraph structure matches the\n# target exactly.\n\n$d = $word.ActiveDocument\n\n$startRange = $d.Content\n$null = $startRange.Find.Execute(\"Pegar a descri\u00e7\u00e3o com o professor\")\n$startParagraph = $startRange.Paragraphs(1)\n\n$endRange = $d.Content\n$null = $endRange.Find.Execute(\"Pegar quais s\u00e3o os\")\n$endParagraph = $endRange.Paragraphs(1)\n\n$targetRange = $d.Range($startParagraph.Range.Start, $endParagraph.Range.End)\n\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:jc w:val=\"center\"/><w:rPr><w:b/><w:sz w:val=\"36\"/><w:szCs w:val=\"36\"/></w:rPr></w:pPr><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p><w:p><w:pPr><w:jc w:val=\"both\"/><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>Problem\u00e1tica:</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val=\"both\"/><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>N\u00e3o \u00e9 incomum a ocorr\u00eancia de erros em aulas de laborat\u00f3rio</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> que envolvem amplificadores operacionais</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\">. Em tais situa\u00e7\u00f5es, o aluno/professor \u00e9 muitas vezes incapaz de percebe se o problema \u00e9 do </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>Amp</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>-OP utilizado</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\">, que acaba por estragar com facilidade durante as aulas pr\u00e1ticas, </w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\">ou da montagem do circuito feito pelo aluno. </w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>Na tentativa de compreender onde est\u00e1 o erro, gasta-se muito do tempo da aula.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val=\"both\"/><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/></w:rPr><w:t>Projeto:</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val=\"both\"/><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>Este projeto visa desenvolver e montar um hardware capaz de</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> testar</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\">a </w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>funcionalidade</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> de um amplificador operacional, para a facilita\u00e7\u00e3o de aulas pr\u00e1ticas</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\">. </w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>Ele s</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>er\u00e1 feito para funcionar com os modelos utilizados no IFSC</w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> \u2013 campus Florian\u00f3polis: </w:t></w:r><w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>LM324, TL082.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$targetRange.InsertXML($ooxml)\n"}
